$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The inventory date value (with its date-formatted style) that used to live
# in E2 actually belongs in F2 - shift it over one column to the right.
$ws.Range("E2").Copy($ws.Range("F2"))
$ws.Range("E2").Clear()

# Reflect the new active cell/selection in the sheet view.
$ws.Range("F2").Select()
